$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.903.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.699.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '2.43'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +29.16%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '229.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '654.50'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.441'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.15'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +8.83%  '
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.697.21'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.210'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000302'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.391.15'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.508.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.695.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.547'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '534.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.252'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +45.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '119.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +19.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000212'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.899.11'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.39'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.01'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.187'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '33.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.614'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("B38").Value = 'Binance-PegBSC-USD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '611.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.507'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +17.78%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0516'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +14.68%  '
$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.164'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.01%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.970'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.68%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
